$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.612.12"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.364.03"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "590.61"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +6.18%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "187.05"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("E9").Value = "  +2.02%  "
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.47"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.55%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.908.96"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.62"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("B15").Value = "BitcoinCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "633.96"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +6.83%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "67.625.59"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.373.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "18.02"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.11"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.24%  "
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.99"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.12"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "99.42"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("E26").Value = "  +4.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.68"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "32.47"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.52%  "
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.98"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "604.81"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.56%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.979.73"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +7.33%  "
$ws.Range("E34").Value = "  +1.49%  "
$ws.Range("E35").Value = "  +1.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.04"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.82"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.60%  "
$ws.Range("E39").Value = "  +4.58%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.24"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0705"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.66%  "
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("E47").Value = "  +1.17%  "
$ws.Range("E48").Value = "  +0.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +11.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.72"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -21.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "129.67"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +4.72%  "
